$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")
if (-not $ws) { $ws = $wb.ActiveSheet }

# Update the "last updated" timestamp title in A1
$ws.Range("A1").Value = "Datos actualizados a 12 de Abril de 2020 a las 18:22"

# Swap Senegal / Estado de Palestina order (rows 108 and 109)
$ws.Range("A108").Value = "Estado de Palestina"
$ws.Range("A109").Value = "Senegal"

# Row 4 - Estados Unidos
$ws.Range("B4").Value = 537210
$ws.Range("C4").Value = 4331
$ws.Range("D4").Value = 31067
$ws.Range("E4").Value = 484715
$ws.Range("F4").Value = 11665
$ws.Range("G4").Value = 851
$ws.Range("H4").Value = 21428

# Row 6 - Italia
$ws.Range("B6").Value = 156363
$ws.Range("C6").Value = 4092
$ws.Range("D6").Value = 34211
$ws.Range("E6").Value = 102253
$ws.Range("F6").Value = 3343
$ws.Range("G6").Value = 431
$ws.Range("H6").Value = 19899

# Row 24 - India
$ws.Range("D24").Value = 1061
$ws.Range("E24").Value = 7780
$ws.Range("G24").Value = 37
$ws.Range("H24").Value = 325

# Row 27 - Chile
$ws.Range("E27").Value = 5074
$ws.Range("F27").Value = 387
$ws.Range("G27").Value = 7
$ws.Range("H27").Value = 80

# Row 31 - Noruega
$ws.Range("E31").Value = 6327
$ws.Range("G31").Value = 7
$ws.Range("H31").Value = 126

# Row 33 - Rumania
$ws.Range("E33").Value = 5138
$ws.Range("G33").Value = 19
$ws.Range("H33").Value = 310

# Row 35 - Chequia
$ws.Range("B35").Value = 5952
$ws.Range("C35").Value = 50
$ws.Range("D35").Value = 464
$ws.Range("E35").Value = 5350
$ws.Range("G35").Value = 9
$ws.Range("H35").Value = 138

# Row 53 - Singapur
$ws.Range("B53").Value = 2532
$ws.Range("C53").Value = 233
$ws.Range("D53").Value = 560
$ws.Range("E53").Value = 1964

# Row 55 - Grecia
$ws.Range("F55").Value = 76

# Row 58 - Argelia
$ws.Range("B58").Value = 1914
$ws.Range("C58").Value = 89
$ws.Range("D58").Value = 591
$ws.Range("E58").Value = 1030
$ws.Range("G58").Value = 18
$ws.Range("H58").Value = 293

# Row 74 - Bosnia y Herzegovina
$ws.Range("B74").Value = 1004
$ws.Range("C74").Value = 58
$ws.Range("E74").Value = 772
$ws.Range("G74").Value = 2
$ws.Range("H74").Value = 39

# Row 103 - Kirguistan
$ws.Range("D103").Value = 54
$ws.Range("E103").Value = 318

# Row 108 - now Estado de Palestina
$ws.Range("B108").Value = 288
$ws.Range("C108").Value = 20
$ws.Range("D108").Value = 58
$ws.Range("E108").Value = 228
$ws.Range("F108").Value = 0

# Row 109 - now Senegal
$ws.Range("B109").Value = 280
$ws.Range("C109").Value = 2
$ws.Range("D109").Value = 171
$ws.Range("E109").Value = 107
$ws.Range("F109").Value = 1
